$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "阳光电源"
$ws.Range("B2").Value = "工业富联"
$ws.Range("C2").Value = "工业富联"

$ws.Range("A3").Value = "工业富联"
$ws.Range("B3").Value = "山子高科"
$ws.Range("C3").Value = "新易盛"

$ws.Range("A4").Value = "平潭发展"
$ws.Range("B4").Value = "阳光电源"
$ws.Range("C4").Value = "江波龙"

$ws.Range("A5").Value = "山子高科"
$ws.Range("B5").Value = "隆基绿能"
$ws.Range("C5").Value = "山子高科"

$ws.Range("A6").Value = "隆基绿能"
$ws.Range("B6").Value = "平潭发展"
$ws.Range("C6").Value = "阳光电源"

$ws.Range("A7").Value = "三花智控"
$ws.Range("B7").Value = "三花智控"
$ws.Range("C7").Value = "平潭发展"

$ws.Range("A8").Value = "中钨高新"
$ws.Range("B8").Value = "多氟多"
$ws.Range("C8").Value = "隆基绿能"

$ws.Range("A9").Value = "先导智能"
$ws.Range("B9").Value = "先导智能"
$ws.Range("C9").Value = "先导智能"

$ws.Range("A10").Value = "神州信息"
$ws.Range("B10").Value = "方正科技"
$ws.Range("C10").Value = "贵州茅台"

$ws.Range("A11").Value = "上海电力"
$ws.Range("B11").Value = "东方财富"
$ws.Range("C11").Value = "上海电力"

$ws.Range("A12").Value = "安泰科技"
$ws.Range("B12").Value = "上海电力"
$ws.Range("C12").Value = "三花智控"

$ws.Range("A13").Value = "铜冠铜箔"
$ws.Range("B13").Value = "特变电工"
$ws.Range("C13").Value = "安泰科技"

$ws.Range("A14").Value = "多氟多"
$ws.Range("B14").Value = "中钨高新"
$ws.Range("C14").Value = "首开股份"

$ws.Range("A15").Value = "科大国创"
$ws.Range("B15").Value = "安泰科技"
$ws.Range("C15").Value = "青岛双星"

$ws.Range("A16").Value = "特变电工"
$ws.Range("B16").Value = "海峡创新"
$ws.Range("C16").Value = "中钨高新"

$ws.Range("A17").Value = "中兴通讯"
$ws.Range("B17").Value = "中国核建"
$ws.Range("C17").Value = "方大炭素"

$ws.Range("A18").Value = "胜宏科技"
$ws.Range("B18").Value = "神州信息"
$ws.Range("C18").Value = "中兴通讯"

$ws.Range("A19").Value = "方大炭素"
$ws.Range("B19").Value = "中兴通讯"
$ws.Range("C19").Value = "天融信"

$ws.Range("A20").Value = "海峡创新"
$ws.Range("B20").Value = "方大炭素"
$ws.Range("C20").Value = "盈新发展"

$ws.Range("A21").Value = "方正科技"
$ws.Range("B21").Value = "铜冠铜箔"
$ws.Range("C21").Value = "合锻智能"
